$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 9 new data rows (rows 22-30), continuing the existing table pattern:
# regcntr_id, usr_id, lang_code="eng", is_active=TRUE, cr_by="superadmin", cr_dtimes="now()"
for ($i = 0; $i -lt 9; $i++) {
    $row = 22 + $i
    $ws.Cells.Item($row, 1).Value = 10002 + $i
    $ws.Cells.Item($row, 2).Value = 110021 + $i
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
}

# Mirror the leftover UI selection state (rows below the new data, through the
# end of the sheet) that Excel persists in the sheetView when the file is saved.
[void]$ws.Rows("31:1048576").Select()

# Set the printable page orientation to portrait (as recorded in pageSetup).
$ws.PageSetup.Orientation = 1

Write-Host "done"
